# Auto-generated Excel COM-interop script
# Applies numeric value updates to market-price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 537.36365
$ws.Range("I33").Value = 523.5185
$ws.Range("K33").Value = 523.5185
$ws.Range("M33").Value = -294.5185
$ws.Range("H41").Value = 8333604
$ws.Range("I41").Value = 15625209
$ws.Range("J41").Value = 341
$ws.Range("K41").Value = 15625209
$ws.Range("L41").Value = 341
$ws.Range("M41").Value = -15624769
$ws.Range("N41").Value = -1221
$ws.Range("H82").Value = 999
$ws.Range("I82").Value = 999
$ws.Range("K82").Value = 2997
$ws.Range("M82").Value = -2591
$ws.Range("H85").Value = 999
$ws.Range("I85").Value = 999
$ws.Range("K85").Value = 2997
$ws.Range("M85").Value = -1593
$ws.Range("H88").Value = 15189281
$ws.Range("I88").Value = 55558324
$ws.Range("J88").Value = 50890.062
$ws.Range("K88").Value = 55558324
$ws.Range("L88").Value = 50890.062
$ws.Range("M88").Value = -55557918
$ws.Range("N88").Value = -51702.062
$ws.Range("H91").Value = 15189281
$ws.Range("I91").Value = 55558324
$ws.Range("J91").Value = 50890.062
$ws.Range("K91").Value = 55558324
$ws.Range("L91").Value = 50890.062
$ws.Range("M91").Value = -55556920
$ws.Range("N91").Value = -53698.062
$ws.Range("H113").Value = 55563504
$ws.Range("I113").Value = 27779882
$ws.Range("J113").Value = 71439860
$ws.Range("K113").Value = 27779882
$ws.Range("L113").Value = 71439860
$ws.Range("M113").Value = -27776628
$ws.Range("N113").Value = -71446368

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2123947.8
$ws.Range("I32").Value = 2455673.2
$ws.Range("K32").Value = 2455673.2
$ws.Range("M32").Value = -2455386.2
$ws.Range("H61").Value = 4945.797
$ws.Range("I61").Value = 2468.2766
$ws.Range("J61").Value = 11795.412
$ws.Range("K61").Value = 2468.2766
$ws.Range("L61").Value = 11795.412
$ws.Range("M61").Value = -2256.2766
$ws.Range("N61").Value = -12219.412
$ws.Range("H97").Value = 8350243
$ws.Range("J97").Value = 27831116
$ws.Range("L97").Value = 27831116
$ws.Range("N97").Value = -27832108
$ws.Range("H132").Value = 4401.288
$ws.Range("I132").Value = 3479.4792
$ws.Range("J132").Value = 6859.4443
$ws.Range("K132").Value = 10438.4376
$ws.Range("L132").Value = 20578.3329
$ws.Range("M132").Value = -7908.437600000001
$ws.Range("N132").Value = -25638.3329
$ws.Range("H136").Value = 4945.797
$ws.Range("I136").Value = 2468.2766
$ws.Range("J136").Value = 11795.412
$ws.Range("K136").Value = 7404.8298
$ws.Range("L136").Value = 35386.236
$ws.Range("M136").Value = -4854.8298
$ws.Range("N136").Value = -40486.236
$ws.Range("H138").Value = 73800
$ws.Range("J138").Value = 73800
$ws.Range("L138").Value = 73800
$ws.Range("N138").Value = -84080

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1886
$ws.Range("I94").Value = 815.6
$ws.Range("K94").Value = 815.6
$ws.Range("M94").Value = -364.6
$ws.Range("H105").Value = 47581.734
$ws.Range("I105").Value = 61086.04
$ws.Range("J105").Value = 3692.75
$ws.Range("K105").Value = 61086.04
$ws.Range("L105").Value = 3692.75
$ws.Range("M105").Value = -59339.04
$ws.Range("N105").Value = -7186.75
$ws.Range("H134").Value = 4812.0566
$ws.Range("I134").Value = 1835.7812
$ws.Range("K134").Value = 5507.3436
$ws.Range("M134").Value = -2972.3436

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8872.342000000001
$ws.Range("I31").Value = 3180.75
$ws.Range("K31").Value = 3180.75
$ws.Range("M31").Value = -2885.75
$ws.Range("H34").Value = 8872.342000000001
$ws.Range("I34").Value = 3180.75
$ws.Range("K34").Value = 3180.75
$ws.Range("M34").Value = -2978.75
$ws.Range("H58").Value = 21748572
$ws.Range("I58").Value = 83335816
$ws.Range("J58").Value = 11896.412
$ws.Range("K58").Value = 83335816
$ws.Range("L58").Value = 11896.412
$ws.Range("M58").Value = -83335613
$ws.Range("N58").Value = -12302.412
$ws.Range("H64").Value = 61841.43
$ws.Range("J64").Value = 61841.43
$ws.Range("L64").Value = 61841.43
$ws.Range("N64").Value = -62337.43
$ws.Range("H67").Value = 61841.43
$ws.Range("J67").Value = 61841.43
$ws.Range("L67").Value = 61841.43
$ws.Range("N67").Value = -63557.43
$ws.Range("H107").Value = 1665.28
$ws.Range("I107").Value = 1595
$ws.Range("J107").Value = 1720.5
$ws.Range("K107").Value = 1595
$ws.Range("L107").Value = 1720.5
$ws.Range("M107").Value = 325
$ws.Range("N107").Value = -5560.5
$ws.Range("H132").Value = 4708.4653
$ws.Range("I132").Value = 1505.3846
$ws.Range("K132").Value = 4516.1538
$ws.Range("M132").Value = -1986.1538
$ws.Range("H136").Value = 21748572
$ws.Range("I136").Value = 83335816
$ws.Range("J136").Value = 11896.412
$ws.Range("K136").Value = 250007448
$ws.Range("L136").Value = 35689.236
$ws.Range("M136").Value = -250004898
$ws.Range("N136").Value = -40789.236

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3125483.5
$ws.Range("I12").Value = 781.3333
$ws.Range("K12").Value = 2343.9999
$ws.Range("M12").Value = -2170.9999
$ws.Range("H63").Value = 4006
$ws.Range("I63").Value = 4006
$ws.Range("K63").Value = 12018
$ws.Range("M63").Value = -11269
$ws.Range("H66").Value = 4006
$ws.Range("I66").Value = 4006
$ws.Range("K66").Value = 36054
$ws.Range("M66").Value = -32310
$ws.Range("H107").Value = 12500294
$ws.Range("I107").Value = 358.5
$ws.Range("J107").Value = 16666939
$ws.Range("K107").Value = 1075.5
$ws.Range("L107").Value = 50000817
$ws.Range("M107").Value = 844.5
$ws.Range("N107").Value = -50004657
$ws.Range("H114").Value = 622.9
$ws.Range("I114").Value = 246.6
$ws.Range("J114").Value = 999.2
$ws.Range("K114").Value = 739.8
$ws.Range("L114").Value = 2997.6
$ws.Range("M114").Value = 2514.2
$ws.Range("N114").Value = -9505.6

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 322.85715
$ws.Range("I17").Value = 276.66666
$ws.Range("K17").Value = 276.66666
$ws.Range("M17").Value = -108.66666
$ws.Range("H97").Value = 283.66666
$ws.Range("I97").Value = 294.0909
$ws.Range("K97").Value = 294.0909
$ws.Range("M97").Value = 201.9091

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5461.048
$ws.Range("I40").Value = 3897.9092
$ws.Range("K40").Value = 3897.9092
$ws.Range("M40").Value = -3761.9092
$ws.Range("H46").Value = 2133.625
$ws.Range("I46").Value = 484
$ws.Range("J46").Value = 3123.4
$ws.Range("K46").Value = 484
$ws.Range("L46").Value = 3123.4
$ws.Range("M46").Value = -296
$ws.Range("N46").Value = -3499.4
$ws.Range("H82").Value = 705530.1
$ws.Range("J82").Value = 2266.5
$ws.Range("L82").Value = 2266.5
$ws.Range("N82").Value = -2988.5
$ws.Range("H85").Value = 705530.1
$ws.Range("J85").Value = 2266.5
$ws.Range("L85").Value = 2266.5
$ws.Range("N85").Value = -4762.5
$ws.Range("H107").Value = 1929.5
$ws.Range("I107").Value = 1929.5
$ws.Range("K107").Value = 1929.5
$ws.Range("M107").Value = -9.5
$ws.Range("H132").Value = 10422798
$ws.Range("I132").Value = 21741144
$ws.Range("J132").Value = 9919.719999999999
$ws.Range("K132").Value = 65223432
$ws.Range("L132").Value = 29759.16
$ws.Range("M132").Value = -65220902
$ws.Range("N132").Value = -34819.16

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2017.6
$ws.Range("J62").Value = 1772
$ws.Range("L62").Value = 1772
$ws.Range("N62").Value = -3020
$ws.Range("H65").Value = 2017.6
$ws.Range("J65").Value = 1772
$ws.Range("L65").Value = 8860
$ws.Range("N65").Value = -15100
$ws.Range("H96").Value = 4277.7144
$ws.Range("I96").Value = 4277.7144
$ws.Range("K96").Value = 4277.7144
$ws.Range("M96").Value = -2904.7144
$ws.Range("H122").Value = 4376.577
$ws.Range("I122").Value = 3717.6365
$ws.Range("J122").Value = 5521.0527
$ws.Range("K122").Value = 11152.9095
$ws.Range("L122").Value = 16563.1581
$ws.Range("M122").Value = -8702.9095
$ws.Range("N122").Value = -21463.1581
$ws.Range("H132").Value = 11120538
$ws.Range("I132").Value = 14709850
$ws.Range("J132").Value = 26300.455
$ws.Range("K132").Value = 44129550
$ws.Range("L132").Value = 78901.36500000001
$ws.Range("M132").Value = -44127020
$ws.Range("N132").Value = -83961.36500000001
